# Add new date columns (JP, JQ, JR) with one more week of mobility data,
# mirroring the style (date number format) of the existing date header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: new date headers (stored as Excel serial dates, formatted like the
# existing header cells in row 1, e.g. JO1). Copy JO1's formatting (date
# number format) into the new cells first so they share the same style.
$ws.Range("JO1").Copy($ws.Range("JP1:JR1"))
$ws.Range("JP1").Value = 44112
$ws.Range("JQ1").Value = 44113
$ws.Range("JR1").Value = 44114

# Row 2 new data
$ws.Range("JP2").Value = 53.56
$ws.Range("JQ2").Value = 61.5
$ws.Range("JR2").Value = 65.53

# Row 3 new data
$ws.Range("JP3").Value = 35.01
$ws.Range("JQ3").Value = 39.29
$ws.Range("JR3").Value = 36.33

# Row 4 new data
$ws.Range("JP4").Value = 59.49
$ws.Range("JQ4").Value = 65.76
$ws.Range("JR4").Value = 60.44

# Row 5 new data
$ws.Range("JP5").Value = 61.72
$ws.Range("JQ5").Value = 65.56
$ws.Range("JR5").Value = 61.4

# Match the author's final selection state recorded in the workbook.
$ws.Range("JT27").Select()
